$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new SMS/notification entry was captured just after the one currently in
# row 31 ("transfer" @ 2024-09-05 16:28:38). Insert a fresh row above it so
# every existing entry (rows 31-75) shifts down by one (to rows 32-76), then
# populate the new row 31 with the latest entry.
$ws.Rows("31:31").Insert()

$ws.Range("R31").Value = "share anyone axis"
$ws.Range("S31").Value = "2024-09-05 16:31:34"
